$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data, shifting everything down.
$ws.Rows("1:1").Insert()

# Populate the new header row.
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Suscripcion"

# Restore the selection to the default top-left cell.
$ws.Range("A1").Select() | Out-Null
